# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" worksheet right after "总计" (i.e. before the
#    existing "2022-Q3" sheet), populated with the new quarter's fund table.
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at the
#    top of the data (row 2) and shift the remaining quarters down by one
#    row, re-numbering the index column (A) sequentially.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: insert & populate the new "2022-Q4" worksheet
# ---------------------------------------------------------------------

$summary   = $wb.Worksheets.Item(1)        # "总计"
$q3Sheet   = $wb.Worksheets.Item(2)        # currently "2022-Q3"

$q4Sheet = $wb.Worksheets.Add($q3Sheet)    # inserted directly before "2022-Q3"
$q4Sheet.Name = "2022-Q4"

# NOTE: after Add() shifts everything right by one, re-resolve "2022-Q3" by
# its new index (3) instead of reusing the pre-insert $q3Sheet reference,
# which gets reseated to point at the freshly-added sheet.
$q3Sheet = $wb.Worksheets.Item(3)

# Clone the cell formatting (header style s=2 across B1:H1 and A2:A8) from
# the existing "2022-Q3" sheet so the new sheet matches the workbook's
# established look without hand-rolling style indices.
$q3Sheet.Range("A1:H8").Copy()
$q4Sheet.Range("A1:H8").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$q4Sheet.Cells.Item(1, 2).Value = "基金代码"
$q4Sheet.Cells.Item(1, 3).Value = "基金名称"
$q4Sheet.Cells.Item(1, 4).Value = "基金规模"
$q4Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q4Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q4Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4Sheet.Cells.Item(1, 8).Value = "仓位排名"

# Data rows. Columns C-G (name/scale/position/weight/value) are text in the
# source data (e.g. "011243" keeps its leading zero, "3.32" stays text) -
# only column H (rank) is a genuine number.
$q4Data = @(
    @("320011", "诺安中小盘精选混合",             "3.32", "84.80", "3.74", "0.1242", 8),
    @("320015", "诺安行业轮动混合",                 "1.21", "86.99", "3.31", "0.0401", 8),
    @("011243", "万家惠裕回报6个月持有期混合A", "1.28", "29.05", "1.74", "0.0223", 1),
    @("519615", "银河君尚灵活配置混合I",         "1.83", "38.98", "1.07", "0.0196", 2),
    @("519613", "银河君尚灵活配置混合A",         "1.17", "38.98", "1.07", "0.0125", 2),
    @("519614", "银河君尚灵活配置混合C",         "0.16", "38.98", "1.07", "0.0017", 2),
    @("011244", "万家惠裕回报6个月持有期混合C", "0.10", "29.05", "1.74", "0.0017", 1)
)

# Force text storage for B:G (so values like "011243"/"3.32" are not
# re-interpreted as numbers and don't lose leading zeros / precision),
# then restore the default (unstyled) look to match the source workbook.
$q4TextRange = $q4Sheet.Range("B2:G8")
$q4TextRange.NumberFormat = "@"

$r = 2
$idx = 0
foreach ($row in $q4Data) {
    $q4Sheet.Cells.Item($r, 1).Value = $idx
    $q4Sheet.Cells.Item($r, 2).Value = $row[0]
    $q4Sheet.Cells.Item($r, 3).Value = $row[1]
    $q4Sheet.Cells.Item($r, 4).Value = $row[2]
    $q4Sheet.Cells.Item($r, 5).Value = $row[3]
    $q4Sheet.Cells.Item($r, 6).Value = $row[4]
    $q4Sheet.Cells.Item($r, 7).Value = $row[5]
    $q4Sheet.Cells.Item($r, 8).Value = $row[6]
    $r++
    $idx++
}

$q4TextRange.Style = "Normal"

# ---------------------------------------------------------------------
# Part 2: update the "总计" summary sheet with the new quarter
# ---------------------------------------------------------------------

$summary.Rows.Item(2).Insert()

# New row 2 needs the A-column index style (s=2) that the rest of the
# index column uses; clone it from row 3 (the old row 2, shifted down).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)      # xlPasteFormats
$summary.Range("B2:D2").Style = "Normal"

$summaryData = @(
    @("2022-Q4", 7,  "0.22"),
    @("2022-Q3", 13, "0.21"),
    @("2022-Q2", 6,  "0.31"),
    @("2022-Q1", 4,  "0.15"),
    @("2021-Q4", 7,  "0.26"),
    @("2021-Q3", 1,  "0.07000000000000001"),
    @("2021-Q1", 18, "0.5")
)

$r = 2
$idx = 0
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = $idx
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = [double]$row[2]
    $r++
    $idx++
}
